$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.592.52"
$ws.Range("E2").Value = "  -7.31%  "

$ws.Range("D3").Value = "2.172.24"
$ws.Range("E3").Value = "  -7.76%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'238.30"
$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").Value = "'0.615"
$ws.Range("E6").Value = "  -7.99%  "

$ws.Range("D7").Value = "'69.24"
$ws.Range("E7").Value = "  -5.52%  "

$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("E9").Value = "  -11.86%  "

$ws.Range("D10").Value = "'36.25"
$ws.Range("E10").Value = "  +3.70%  "

$ws.Range("D11").Value = "'57.36"
$ws.Range("E11").Value = "  -5.72%  "

$ws.Range("D12").Value = "'0.0931"
$ws.Range("E12").Value = "  -8.85%  "

$ws.Range("E13").Value = "  -4.74%  "

$ws.Range("D14").Value = "'6.47"
$ws.Range("E14").Value = "  -10.14%  "

$ws.Range("D15").Value = "2.501.18"
$ws.Range("E15").Value = "  -7.63%  "

$ws.Range("D16").Value = "'14.49"
$ws.Range("E16").Value = "  -10.44%  "

$ws.Range("E17").Value = "  -9.68%  "

$ws.Range("D18").Value = "2.174.23"
$ws.Range("E18").Value = "  -7.66%  "

$ws.Range("D19").Value = "40.616.57"
$ws.Range("E19").Value = "  -7.19%  "

$ws.Range("D20").Value = "0.0₃0930"
$ws.Range("E20").Value = "  -9.66%  "

$ws.Range("D21").Value = "'71.71"
$ws.Range("E21").Value = "  -7.64%  "

$ws.Range("E22").Value = "  -8.25%  "

$ws.Range("D23").Value = "'228.38"
$ws.Range("E23").Value = "  -9.01%  "

$ws.Range("E24").Value = "  +5.66%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "'3.59"
$ws.Range("E26").Value = "  -4.93%  "

$ws.Range("E27").Value = "  -4.64%  "

$ws.Range("E28").Value = "  -5.25%  "

$ws.Range("E29").Value = "  -8.56%  "

$ws.Range("D30").Value = "'168.16"
$ws.Range("E30").Value = "  -4.25%  "

$ws.Range("D31").Value = "'19.98"
$ws.Range("E31").Value = "  -10.28%  "

$ws.Range("E32").Value = "  -10.05%  "

$ws.Range("E33").Value = "  -8.53%  "

$ws.Range("E34").Value = "  -7.58%  "

$ws.Range("D35").Value = "'5.06"
$ws.Range("E35").Value = "  -5.32%  "

$ws.Range("E36").Value = "  -10.17%  "

$ws.Range("D37").Value = "'3.74"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").Value = "'22.81"
$ws.Range("E38").Value = "  +14.19%  "

$ws.Range("E39").Value = "  -7.86%  "

$ws.Range("D40").Value = "'0.0263"
$ws.Range("E40").Value = "  -4.65%  "

$ws.Range("D41").Value = "'5.77"
$ws.Range("E41").Value = "  -12.45%  "

$ws.Range("D42").Value = "'63.20"
$ws.Range("E42").Value = "  -2.59%  "

$ws.Range("D43").Value = "'4.76"
$ws.Range("E43").Value = "  -14.00%  "

$ws.Range("D44").Value = "'8.56"
$ws.Range("E44").Value = "  -5.08%  "

$ws.Range("E45").Value = "  -6.73%  "

$ws.Range("E46").Value = "  +0.24%  "

$ws.Range("E47").Value = "  -8.29%  "

$ws.Range("D48").Value = "'4.42"
$ws.Range("E48").Value = "  +1.48%  "

$ws.Range("D49").Value = "'10.18"
$ws.Range("E49").Value = "  +6.88%  "

$ws.Range("E50").Value = "  -6.58%  "

$ws.Range("E51").Value = "  -7.00%  "
